$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 610, shifting existing rows 610:657 down to 611:658
$ws.Rows(610).Insert()

# Populate the newly inserted row 610 with the new record's data
$ws.Cells.Item(610, 1).Value = 3
$ws.Cells.Item(610, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(610, 3).Value = "Coquimbo"
$ws.Cells.Item(610, 4).Value = 45013
$ws.Cells.Item(610, 5).Value = 5
$ws.Cells.Item(610, 6).Value = 100112021
$ws.Cells.Item(610, 7).Value = "Ají"
$ws.Cells.Item(610, 8).Value = "Inferno"
$ws.Cells.Item(610, 9).Value = "Primera"
$ws.Cells.Item(610, 10).Value = 65
$ws.Cells.Item(610, 11).Value = 22000
$ws.Cells.Item(610, 12).Value = 23000
$ws.Cells.Item(610, 13).Value = 22538
$ws.Cells.Item(610, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(610, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(610, 16).Value = 1503
$ws.Cells.Item(610, 17).Value = 15
$ws.Cells.Item(610, 18).Value = "Hortaliza"
